$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update event names (rows 2-9 in column A) due to Daylight Savings Time change.
$ws.Range("A2").Value = "Defend the Vault 1"
$ws.Range("A3").Value = "Defend the Vault 2"
$ws.Range("A4").Value = "Battlegrounds 1"
$ws.Range("A5").Value = "Battlegrounds 2"
$ws.Range("A6").Value = "Battlegrounds 3"
$ws.Range("A7").Value = "Battlegrounds 4"
$ws.Range("A8").Value = "Corvus Expedition"
$ws.Range("A9").Value = "Rite of Exile"

# Move the selection to A9, matching the saved workbook state.
$ws.Range("A9").Select()
